# Auto-generated edit script: update "想去人数" (F column) counts
# and sold-out status (G column) for two rows, per commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1509
$ws.Range("F3").Value = 1473
$ws.Range("F6").Value = 742
$ws.Range("F7").Value = 45
$ws.Range("F8").Value = 674
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 1396
$ws.Range("F12").Value = 36668
$ws.Range("F13").Value = 7291
$ws.Range("F14").Value = 119
$ws.Range("F15").Value = 385
$ws.Range("F16").Value = 595
$ws.Range("F17").Value = 459
$ws.Range("F21").Value = 55
$ws.Range("F22").Value = 464
$ws.Range("F23").Value = 127
$ws.Range("F24").Value = 834
$ws.Range("F25").Value = 24
$ws.Range("F26").Value = 329
$ws.Range("F30").Value = 234
$ws.Range("F31").Value = 62
$ws.Range("F32").Value = 755
$ws.Range("F33").Value = 298
$ws.Range("F35").Value = 772
$ws.Range("F38").Value = 819
$ws.Range("F39").Value = 299

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1234
$ws.Range("F9").Value = 243
$ws.Range("F16").Value = 55
$ws.Range("F19").Value = 4311

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1498
$ws.Range("F3").Value = 372

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1498
$ws.Range("F3").Value = 372
$ws.Range("F4").Value = 1234
$ws.Range("F5").Value = 1509
$ws.Range("F7").Value = 1473
$ws.Range("F9").Value = 742
$ws.Range("F10").Value = 45
$ws.Range("F11").Value = 674
$ws.Range("F13").Value = 1396
$ws.Range("F14").Value = 36668
$ws.Range("F17").Value = 243
$ws.Range("F20").Value = 7291
$ws.Range("F21").Value = 385
$ws.Range("F23").Value = 595
$ws.Range("F24").Value = 459
$ws.Range("F28").Value = 55
$ws.Range("F30").Value = 464
$ws.Range("F31").Value = 127
$ws.Range("F32").Value = 834
$ws.Range("F33").Value = 24
$ws.Range("F34").Value = 329
$ws.Range("F38").Value = 234
$ws.Range("F39").Value = 62
$ws.Range("F40").Value = 755
$ws.Range("F41").Value = 55
$ws.Range("F42").Value = 298
$ws.Range("F44").Value = 819
$ws.Range("F45").Value = 299

# Sold-out status update for "萤火虫动漫游戏嘉年华 x KKWORLD2024" row
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G12").Value = "暂时售罄"

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G14").Value = 85
